# Update schedule for fall 2017
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Structural changes
# ------------------------------------------------------------------
# Row 5 ("cm000" / "No class (MLK Jr. Day)") is removed entirely - the
# cm-section below it shifts up by one row.
$ws.Rows("5:5").Delete()

# Two new course-meeting rows (cm019 / cm020) are added right before the
# lab section (which, after the deletion above, starts at row 20).
$ws.Rows("20:21").Insert()

# ------------------------------------------------------------------
# 2. Introduce the brand new label text in the same order the workbook
#    originally introduced them, so sharedStrings ends up in the same
#    append order.
# ------------------------------------------------------------------
$ws.Range("A20").Value = "cm019"
$ws.Range("A21").Value = "cm020"
$ws.Range("D4").Value = "Data transformation"
$ws.Range("D5").Value = "Exploratory data analysis"
$ws.Range("D7").Value = "Data wrangling (more)"
$ws.Range("D21").Value = "Shiny applications"
$ws.Range("D20").Value = "Interactivity in R"
$ws.Range("D8").Value = "Pipes and functions"
$ws.Range("D11").Value = "R Markdown"
$ws.Range("D25").Value = "The shell"
$ws.Range("D28").Value = "AWS"

# ------------------------------------------------------------------
# 3. Re-write the "cm" (course meeting) section - rows 2 through 21
# ------------------------------------------------------------------

$ws.Range("B2").Value = 43003
$ws.Range("D2").Value = "Introduction to computing for the social sciences"

$ws.Range("B3").Formula = "=B2+2"
$ws.Range("D3").Value = "Visualizations and the grammar of graphics"

$ws.Range("B4").Formula = "=B2+7"

$ws.Range("B5").Formula = "=B4+2"

$ws.Range("B6").Formula = "=B4+7"
$ws.Range("D6").Value = "Data wrangling"

$ws.Range("B7").Formula = "=B6+2"

$ws.Range("B8").Formula = "=B6+7"

$ws.Range("B9").Formula = "=B8+2"
$ws.Range("D9").Value = "Vectors and iteration"

$ws.Range("B10").Formula = "=B8+7"
$ws.Range("D10").Value = "Debugging and defensive programming"

$ws.Range("B11").Formula = "=B10+2"

$ws.Range("B12").Formula = "=B10+7"
$ws.Range("D12").Value = "Statistical learning: basics and linear regression"

$ws.Range("B13").Formula = "=B12+2"
$ws.Range("D13").Value = "Statistical learning: classification"

$ws.Range("B14").Formula = "=B12+7"
$ws.Range("D14").Value = "Statistical learning: resampling methods"

$ws.Range("B15").Formula = "=B14+2"
$ws.Range("D15").Value = "Distributed learning"

$ws.Range("B16").Formula = "=B14+7"
$ws.Range("D16").Value = "Getting data from the web: API access"

$ws.Range("B17").Formula = "=B16+2"
$ws.Range("D17").Value = "Getting data from the web: scraping"

$ws.Range("B18").Formula = "=B16+7"
$ws.Range("D18").Value = "Text analysis: fundamentals and sentiment analysis"

$ws.Range("B19").Formula = "=B18+2"
$ws.Range("D19").Value = "Text analysis: classification and topic modeling"

$ws.Range("B20").Formula = "=B18+7"
$ws.Range("C20").Value = $false

$ws.Range("B21").Formula = "=B20+2"
$ws.Range("C21").Value = $true

# ------------------------------------------------------------------
# 4. Re-write the "lab" section - rows 22 through 31
# ------------------------------------------------------------------

$ws.Range("B22").Value = 43005
$ws.Range("D22").Value = "Software setup"

$ws.Range("B23").Formula = "=B22+7"
$ws.Range("D23").Value = "R basics"

$ws.Range("B24").Formula = "=B23+7"
$ws.Range("C24").Value = $true
$ws.Range("D24").Value = "Graphing tips for ``ggplot2`` (and life)"

$ws.Range("B25").Formula = "=B24+7"

$ws.Range("B26").Formula = "=B25+7"

$ws.Range("B27").Formula = "=B26+7"

$ws.Range("B28").Formula = "=B27+7"

$ws.Range("B29").Formula = "=B28+7"

$ws.Range("B30").Formula = "=B29+7"

$ws.Range("A31").Value = "lab10"
$ws.Range("C31").Value = $false
$ws.Range("B31").Formula = "=B30+7"
